$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E6").Value = "SI"
$ws.Range("E7").Value = "SI"
$ws.Range("E8").Value = "SI"

$ws.Range("E9").Select()
